{"js": "// Update the NSW Covid report figures to the 2022-01-23 output values.\n// Each pair is [old exact text, new exact text] as it appears in a w:t run.\nconst replacements = [\n  // Deaths table\n  [\"1154.0\", \"1154\"],\n  [\"39.0\", \"38\"],\n  [\"1371.0\", \"1371\"],\n  [\"33.0\", \"32\"],\n  [\"1655.0\", \"1655\"],\n  [\"2.0\", \"2\"],\n  // Death rate prediction paragraph\n  [\"1669.0\", \"1668\"],\n  // Number of confirmed infections table\n  [\"994703.0\", \"994883\"],\n  [\"17393.0\", \"17366\"],\n  [\"1058530.0\", \"1058601\"],\n  [\"6909.0\", \"6896\"],\n  [\"1084614.0\", \"1084635\"],\n  [\"2099.0\", \"2094\"],\n  [\"1094864.0\", \"1094865\"],\n  [\"124.0\", \"123\"],\n  // Final number of infections prediction paragraph\n  [\"1095497.0\", \"1095496\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the NSW Covid report figures to the 2022-01-23 output values.\n$d = $word.ActiveDocument\n\n# Each pair is the exact old text of a w:t run and its replacement text.\n$pairs = @(\n    # Deaths table\n    @{old = \"1154.0\"; new = \"1154\"},\n    @{old = \"39.0\";   new = \"38\"},\n    @{old = \"1371.0\"; new = \"1371\"},\n    @{old = \"33.0\";   new = \"32\"},\n    @{old = \"1655.0\"; new = \"1655\"},\n    @{old = \"2.0\";    new = \"2\"},\n    # Death rate prediction paragraph\n    @{old = \"1669.0\"; new = \"1668\"},\n    # Number of confirmed infections table\n    @{old = \"994703.0\";  new = \"994883\"},\n    @{old = \"17393.0\";   new = \"17366\"},\n    @{old = \"1058530.0\"; new = \"1058601\"},\n    @{old = \"6909.0\";    new = \"6896\"},\n    @{old = \"1084614.0\"; new = \"1084635\"},\n    @{old = \"2099.0\";    new = \"2094\"},\n    @{old = \"1094864.0\"; new = \"1094865\"},\n    @{old = \"124.0\";     new = \"123\"},\n    # Final number of infections prediction paragraph\n    @{old = \"1095497.0\"; new = \"1095496\"}\n)\n\nforeach ($pair in $pairs) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Text = $pair.old\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Could not find text: $($pair.old)\"\n    }\n\n    $rng.Text = $pair.new\n}\n"}
